$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 518.5714
$ws.Range("I4").Value = 497
$ws.Range("K4").Value = 497
$ws.Range("M4").Value = -383
$ws.Range("H98").Value = 2854.3333
$ws.Range("I98").Value = 2964.8333
$ws.Range("K98").Value = 2964.8333
$ws.Range("M98").Value = -1466.8333
$ws.Range("H111").Value = 2323.5
$ws.Range("I111").Value = 2323.5
$ws.Range("K111").Value = 6970.5
$ws.Range("M111").Value = -3903.5
$ws.Range("H122").Value = 2854.3333
$ws.Range("I122").Value = 2964.8333
$ws.Range("K122").Value = 8894.499899999999
$ws.Range("M122").Value = -6444.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2776.2354
$ws.Range("I61").Value = 2831.4614
$ws.Range("K61").Value = 2831.4614
$ws.Range("M61").Value = -2619.4614
$ws.Range("H74").Value = 2331.9546
$ws.Range("I74").Value = 2149.8333
$ws.Range("J74").Value = 2550.5
$ws.Range("K74").Value = 2149.8333
$ws.Range("L74").Value = 2550.5
$ws.Range("M74").Value = -1275.8333
$ws.Range("N74").Value = -4298.5
$ws.Range("H77").Value = 2331.9546
$ws.Range("I77").Value = 2149.8333
$ws.Range("J77").Value = 2550.5
$ws.Range("K77").Value = 10749.1665
$ws.Range("L77").Value = 12752.5
$ws.Range("M77").Value = -6381.166499999999
$ws.Range("N77").Value = -21488.5
$ws.Range("H122").Value = 1974.375
$ws.Range("I122").Value = 1466.3334
$ws.Range("K122").Value = 4399.0002
$ws.Range("M122").Value = -1949.0002
$ws.Range("H136").Value = 2776.2354
$ws.Range("I136").Value = 2831.4614
$ws.Range("K136").Value = 8494.3842
$ws.Range("M136").Value = -5944.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("H31").Value = 5534
$ws.Range("I31").Value = 3538.7144
$ws.Range("J31").Value = 10189.667
$ws.Range("K31").Value = 3538.7144
$ws.Range("L31").Value = 10189.667
$ws.Range("M31").Value = -3243.7144
$ws.Range("N31").Value = -10779.667
$ws.Range("H34").Value = 5534
$ws.Range("I34").Value = 3538.7144
$ws.Range("J34").Value = 10189.667
$ws.Range("K34").Value = 3538.7144
$ws.Range("L34").Value = 10189.667
$ws.Range("M34").Value = -3336.7144
$ws.Range("N34").Value = -10593.667
$ws.Range("H62").Value = 12232.333
$ws.Range("I62").Value = 12232.333
$ws.Range("K62").Value = 12232.333
$ws.Range("M62").Value = -11608.333
$ws.Range("H65").Value = 12232.333
$ws.Range("I65").Value = 12232.333
$ws.Range("K65").Value = 61161.665
$ws.Range("M65").Value = -58041.665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 571.8
$ws.Range("I18").Value = 571.8
$ws.Range("K18").Value = 1715.4
$ws.Range("M18").Value = -1546.4
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 51117
$ws.Range("J95").Value = 51117
$ws.Range("L95").Value = 51117
$ws.Range("N95").Value = -56609

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4593.231
$ws.Range("I7").Value = 4923
$ws.Range("K7").Value = 4923
$ws.Range("M7").Value = -4811
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H50").Value = 34500
$ws.Range("I50").Value = 7000
$ws.Range("K50").Value = 7000
$ws.Range("M50").Value = -6363
$ws.Range("H101").Value = 4124.75
$ws.Range("J101").Value = 4124.75
$ws.Range("L101").Value = 4124.75
$ws.Range("N101").Value = -10614.75
$ws.Range("H122").Value = 3670.3928
$ws.Range("I122").Value = 3553.95
$ws.Range("J122").Value = 3961.5
$ws.Range("K122").Value = 10661.85
$ws.Range("L122").Value = 11884.5
$ws.Range("M122").Value = -8211.849999999999
$ws.Range("N122").Value = -16784.5
$ws.Range("H126").Value = 4593.231
$ws.Range("I126").Value = 4923
$ws.Range("K126").Value = 14769
$ws.Range("M126").Value = -12299

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 48186.5
$ws.Range("I34").Value = 45582
$ws.Range("K34").Value = 45582
$ws.Range("M34").Value = -45379
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H40").Value = 52995
$ws.Range("I40").Value = 49990
$ws.Range("K40").Value = 49990
$ws.Range("M40").Value = -49841
$ws.Range("H76").Value = 59999.855
$ws.Range("J76").Value = 59999.855
$ws.Range("L76").Value = 59999.855
$ws.Range("N76").Value = -60629.855
$ws.Range("H79").Value = 59999.855
$ws.Range("J79").Value = 59999.855
$ws.Range("L79").Value = 59999.855
$ws.Range("N79").Value = -62183.855
$ws.Range("H104").Value = 15000
$ws.Range("J104").Value = 15000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -21988
$ws.Range("H107").Value = 719
$ws.Range("I107").Value = 1859.5
$ws.Range("J107").Value = 338.83334
$ws.Range("K107").Value = 5578.5
$ws.Range("L107").Value = 1016.50002
$ws.Range("M107").Value = -3658.5
$ws.Range("N107").Value = -4856.50002
$ws.Range("H113").Value = 1423.375
$ws.Range("I113").Value = 850
$ws.Range("J113").Value = 1996.75
$ws.Range("K113").Value = 2550
$ws.Range("L113").Value = 5990.25
$ws.Range("M113").Value = -380
$ws.Range("N113").Value = -10330.25
$ws.Range("H122").Value = 3293.0789
$ws.Range("I122").Value = 3074.2354
$ws.Range("J122").Value = 5153.25
$ws.Range("K122").Value = 9222.706200000001
$ws.Range("L122").Value = 15459.75
$ws.Range("M122").Value = -6772.706200000001
$ws.Range("N122").Value = -20359.75
$ws.Range("H136").Value = 3667.7693
$ws.Range("I136").Value = 3198.7778
$ws.Range("K136").Value = 9596.3334
$ws.Range("M136").Value = -7046.3334
